# Update 'paises' country data and labels to reflect the refreshed dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 28 de Marzo de 2020 a las 13:29'

# Row 27
$ws.Cells.Item(27, 1).Value = 'Dinamarca'
$ws.Cells.Item(27, 2).Value = 2201
$ws.Cells.Item(27, 3).Value = 155
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(27, 5).Value = 2148
$ws.Cells.Item(27, 6).Value = 109
$ws.Cells.Item(27, 8).Value = 52

# Row 28
$ws.Cells.Item(28, 1).Value = 'Irlanda'
$ws.Cells.Item(28, 2).Value = 2121
$ws.Cells.Item(28, 4).Value = 5
$ws.Cells.Item(28, 5).Value = 2094
$ws.Cells.Item(28, 6).Value = 59
$ws.Cells.Item(28, 8).Value = 22

# Row 87
$ws.Cells.Item(87, 5).Value = 156
$ws.Cells.Item(87, 7).Value = 2
$ws.Cells.Item(87, 8).Value = 10

# Row 92
$ws.Cells.Item(92, 1).Value = 'Islas Feroe'
$ws.Cells.Item(92, 2).Value = 155
$ws.Cells.Item(92, 3).Value = 11
$ws.Cells.Item(92, 4).Value = 54
$ws.Cells.Item(92, 5).Value = 101
$ws.Cells.Item(92, 6).Value = 2

# Row 93
$ws.Cells.Item(93, 1).Value = 'Oman'
$ws.Cells.Item(93, 2).Value = 152
$ws.Cells.Item(93, 3).Value = 21
$ws.Cells.Item(93, 4).Value = 23
$ws.Cells.Item(93, 5).Value = 129
$ws.Cells.Item(93, 6).Value = 0

# Row 94
$ws.Cells.Item(94, 1).Value = 'Malta'
$ws.Cells.Item(94, 2).Value = 149
$ws.Cells.Item(94, 3).Value = 10
$ws.Cells.Item(94, 4).Value = 2
$ws.Cells.Item(94, 5).Value = 147
$ws.Cells.Item(94, 6).Value = 1

# Row 95
$ws.Cells.Item(95, 1).Value = 'Reunion'
$ws.Cells.Item(95, 2).Value = 145
$ws.Cells.Item(95, 4).Value = 1
$ws.Cells.Item(95, 5).Value = 144
$ws.Cells.Item(95, 6).Value = 0

# Row 106
$ws.Cells.Item(106, 1).Value = 'Guadalupe'
$ws.Cells.Item(106, 2).Value = 96
$ws.Cells.Item(106, 3).Value = 23
$ws.Cells.Item(106, 4).Value = 17
$ws.Cells.Item(106, 5).Value = 77
$ws.Cells.Item(106, 7).Value = 1
$ws.Cells.Item(106, 8).Value = 2

# Row 107
$ws.Cells.Item(107, 1).Value = 'Honduras'
$ws.Cells.Item(107, 2).Value = 95
$ws.Cells.Item(107, 3).Value = 27
$ws.Cells.Item(107, 4).Value = 3
$ws.Cells.Item(107, 5).Value = 91
$ws.Cells.Item(107, 6).Value = 4
$ws.Cells.Item(107, 8).Value = 1

# Row 108
$ws.Cells.Item(108, 1).Value = 'Mauricio'
$ws.Cells.Item(108, 4).Value = 0
$ws.Cells.Item(108, 5).Value = 92
$ws.Cells.Item(108, 6).Value = 1
$ws.Cells.Item(108, 8).Value = 2

# Row 109
$ws.Cells.Item(109, 1).Value = 'Bielorrusia'
$ws.Cells.Item(109, 2).Value = 94
$ws.Cells.Item(109, 4).Value = 32
$ws.Cells.Item(109, 5).Value = 62
$ws.Cells.Item(109, 6).Value = 2
$ws.Cells.Item(109, 8).Value = 0

# Row 110
$ws.Cells.Item(110, 1).Value = 'Martinica'
$ws.Cells.Item(110, 2).Value = 93
$ws.Cells.Item(110, 4).Value = 0
$ws.Cells.Item(110, 5).Value = 92
$ws.Cells.Item(110, 6).Value = 12
$ws.Cells.Item(110, 8).Value = 1

# Row 111
$ws.Cells.Item(111, 1).Value = 'Camerun'
$ws.Cells.Item(111, 2).Value = 91
$ws.Cells.Item(111, 3).Value = 0
$ws.Cells.Item(111, 4).Value = 2
$ws.Cells.Item(111, 5).Value = 87
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 8).Value = 2

# Row 112
$ws.Cells.Item(112, 1).Value = 'Georgia'
$ws.Cells.Item(112, 2).Value = 85
$ws.Cells.Item(112, 3).Value = 2
$ws.Cells.Item(112, 4).Value = 14
$ws.Cells.Item(112, 5).Value = 71
$ws.Cells.Item(112, 8).Value = 0

# Row 113
$ws.Cells.Item(113, 1).Value = 'Montenegro'
$ws.Cells.Item(113, 2).Value = 82
$ws.Cells.Item(113, 3).Value = 0
$ws.Cells.Item(113, 4).Value = 0
$ws.Cells.Item(113, 5).Value = 81
$ws.Cells.Item(113, 6).Value = 1

# Row 114
$ws.Cells.Item(114, 1).Value = 'Nigeria'
$ws.Cells.Item(114, 2).Value = 81
$ws.Cells.Item(114, 3).Value = 11
$ws.Cells.Item(114, 4).Value = 3
$ws.Cells.Item(114, 5).Value = 77
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 8).Value = 1

# Row 115
$ws.Cells.Item(115, 1).Value = 'Cuba'
$ws.Cells.Item(115, 2).Value = 80
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(115, 4).Value = 4
$ws.Cells.Item(115, 6).Value = 2
$ws.Cells.Item(115, 8).Value = 2

# Row 116
$ws.Cells.Item(116, 1).Value = 'Bolivia'
$ws.Cells.Item(116, 2).Value = 74
$ws.Cells.Item(116, 3).Value = 13
$ws.Cells.Item(116, 5).Value = 74
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 8).Value = 0

# Row 128
$ws.Cells.Item(128, 1).Value = 'Kenia'
$ws.Cells.Item(128, 2).Value = 38
$ws.Cells.Item(128, 3).Value = 7
$ws.Cells.Item(128, 4).Value = 1
$ws.Cells.Item(128, 5).Value = 36
$ws.Cells.Item(128, 6).Value = 2
$ws.Cells.Item(128, 8).Value = 1

# Row 129
$ws.Cells.Item(129, 1).Value = 'Macao'
$ws.Cells.Item(129, 2).Value = 34
$ws.Cells.Item(129, 4).Value = 10
$ws.Cells.Item(129, 5).Value = 24

# Row 130
$ws.Cells.Item(130, 1).Value = 'Aruba'
$ws.Cells.Item(130, 2).Value = 33
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 1

# Row 131
$ws.Cells.Item(131, 1).Value = 'Isla de Man'
$ws.Cells.Item(131, 3).Value = 3
$ws.Cells.Item(131, 5).Value = 32
$ws.Cells.Item(131, 8).Value = 0

# Row 132
$ws.Cells.Item(132, 1).Value = 'Guam'
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(132, 4).Value = 0
$ws.Cells.Item(132, 5).Value = 31
$ws.Cells.Item(132, 6).Value = 0

# Row 133
$ws.Cells.Item(133, 1).Value = 'Guatemala'
$ws.Cells.Item(133, 2).Value = 32
$ws.Cells.Item(133, 3).Value = 4
$ws.Cells.Item(133, 4).Value = 4
$ws.Cells.Item(133, 5).Value = 27
$ws.Cells.Item(133, 6).Value = 1

# Row 143
$ws.Cells.Item(143, 1).Value = 'Mali'
$ws.Cells.Item(143, 2).Value = 18
$ws.Cells.Item(143, 3).Value = 7
$ws.Cells.Item(143, 5).Value = 18

# Row 144
$ws.Cells.Item(144, 1).Value = 'Islas Virgenes de los Estados Unidos'
$ws.Cells.Item(144, 4).Value = 0
$ws.Cells.Item(144, 5).Value = 17

# Row 145
$ws.Cells.Item(145, 1).Value = 'Bermudas'
$ws.Cells.Item(145, 2).Value = 17
$ws.Cells.Item(145, 4).Value = 2
$ws.Cells.Item(145, 5).Value = 15

# Row 146
$ws.Cells.Item(146, 1).Value = 'Etiopia'
$ws.Cells.Item(146, 4).Value = 1
$ws.Cells.Item(146, 5).Value = 15

# Row 147
$ws.Cells.Item(147, 1).Value = 'Maldivas'
$ws.Cells.Item(147, 2).Value = 16
$ws.Cells.Item(147, 4).Value = 9
$ws.Cells.Item(147, 5).Value = 7

# Row 148
$ws.Cells.Item(148, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(148, 2).Value = 15
$ws.Cells.Item(148, 4).Value = 0
$ws.Cells.Item(148, 5).Value = 15

# Row 149
$ws.Cells.Item(149, 1).Value = 'Tanzania'
$ws.Cells.Item(149, 2).Value = 13
$ws.Cells.Item(149, 4).Value = 1

# Row 150
$ws.Cells.Item(150, 1).Value = 'Republica de Yibuti'

# Row 151
$ws.Cells.Item(151, 1).Value = 'Guinea Ecuatorial'
$ws.Cells.Item(151, 3).Value = 0

# Row 152
$ws.Cells.Item(152, 1).Value = 'Mongolia'
$ws.Cells.Item(152, 2).Value = 12
$ws.Cells.Item(152, 3).Value = 1
$ws.Cells.Item(152, 5).Value = 12

# Row 153
$ws.Cells.Item(153, 1).Value = 'Dominica'

# Row 154
$ws.Cells.Item(154, 1).Value = 'San Martin (Parte Francesa)'

# Row 159
$ws.Cells.Item(159, 1).Value = 'Guinea'

# Row 160
$ws.Cells.Item(160, 1).Value = 'Haiti'

# Row 161
$ws.Cells.Item(161, 1).Value = 'Birmania'

# Row 162
$ws.Cells.Item(162, 1).Value = 'Surinam'

# Row 166
$ws.Cells.Item(166, 1).Value = 'Mozambique'

# Row 167
$ws.Cells.Item(167, 1).Value = 'Granada'

# Row 168
$ws.Cells.Item(168, 1).Value = 'Antigua y Barbuda'

# Row 169
$ws.Cells.Item(169, 1).Value = 'Seychelles'

# Row 175
$ws.Cells.Item(175, 1).Value = 'San Bartolome'
$ws.Cells.Item(175, 3).Value = 0

# Row 176
$ws.Cells.Item(176, 1).Value = 'Fiyi'

# Row 177
$ws.Cells.Item(177, 1).Value = 'Montserrat'

# Row 179
$ws.Cells.Item(179, 1).Value = 'Mauritania'
$ws.Cells.Item(179, 3).Value = 2

# Row 180
$ws.Cells.Item(180, 1).Value = 'Guyana'
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 8).Value = 1

# Row 181
$ws.Cells.Item(181, 1).Value = 'Nepal'
$ws.Cells.Item(181, 3).Value = 1
$ws.Cells.Item(181, 4).Value = 1
$ws.Cells.Item(181, 8).Value = 0

# Row 182
$ws.Cells.Item(182, 1).Value = 'Sudan'
$ws.Cells.Item(182, 3).Value = 2

# Row 183
$ws.Cells.Item(183, 1).Value = 'Cabo Verde'
$ws.Cells.Item(183, 3).Value = 0

# Row 184
$ws.Cells.Item(184, 1).Value = 'Santa Sede'

# Row 185
$ws.Cells.Item(185, 1).Value = 'Congo'

# Row 186
$ws.Cells.Item(186, 1).Value = 'Angola'

# Row 188
$ws.Cells.Item(188, 1).Value = 'Republica de Africa Central'

# Row 189
$ws.Cells.Item(189, 1).Value = 'San Martin (Parte Holandesa)'

# Row 190
$ws.Cells.Item(190, 1).Value = 'Butan'

# Row 192
$ws.Cells.Item(192, 1).Value = 'Liberia'

# Row 193
$ws.Cells.Item(193, 1).Value = 'Somalia'

# Row 196
$ws.Cells.Item(196, 1).Value = 'Islas Turcas y Caicos'

# Row 197
$ws.Cells.Item(197, 1).Value = 'Anguila'

# Row 198
$ws.Cells.Item(198, 1).Value = 'Guinea-Bisau'

# Row 199
$ws.Cells.Item(199, 1).Value = 'San Cristobal y Nieves'

# Row 200
$ws.Cells.Item(200, 1).Value = 'Islas Virgenes Britanicas'

# Row 201
$ws.Cells.Item(201, 1).Value = 'Belice'

# Row 202
$ws.Cells.Item(202, 1).Value = 'Papua Nueva Guinea'

# Row 204
$ws.Cells.Item(204, 1).Value = 'San Vicente y las Granadinas'
